$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D to text format so numeric-looking price
# strings (e.g. "235.21") are not auto-converted to numbers by Excel.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "30.173.73"
$ws.Range("E2").Value = "  -0.20%  "

$ws.Range("D3").Value = "1.849.81"
$ws.Range("E3").Value = "  -0.78%  "

$ws.Range("E4").Value = "  +0.23%  "

$ws.Range("D5").Value = "235.21"
$ws.Range("E5").Value = "  +0.10%  "

$ws.Range("E6").Value = "  +0.20%  "

$ws.Range("D7").Value = "0.4707"
$ws.Range("E7").Value = "  +0.77%  "

$ws.Range("D8").Value = "0.2895"
$ws.Range("E8").Value = "  +2.24%  "

$ws.Range("D9").Value = "0.06523"
$ws.Range("E9").Value = "  +0.20%  "

$ws.Range("D10").Value = "21.67"
$ws.Range("E10").Value = "  +1.56%  "

$ws.Range("D11").Value = "0.07943"
$ws.Range("E11").Value = "  +1.12%  "

$ws.Range("D12").Value = "97.40"
$ws.Range("E12").Value = "  +0.14%  "

$ws.Range("D13").Value = "1.853.87"
$ws.Range("E13").Value = "  -0.60%  "

$ws.Range("D14").Value = "5.085"
$ws.Range("E14").Value = "  -0.11%  "

$ws.Range("D15").Value = "0.6733"
$ws.Range("E15").Value = "  +0.26%  "

$ws.Range("D16").Value = "265.72"
$ws.Range("E16").Value = "  -4.98%  "

$ws.Range("D17").Value = "30.151.61"
$ws.Range("E17").Value = "  -0.25%  "

$ws.Range("D18").Value = "13.59"
$ws.Range("E18").Value = "  +7.50%  "

$ws.Range("D19").Value = "1.003"
$ws.Range("E19").Value = "  +0.27%  "

$ws.Range("D20").Value = "0.000007542"
$ws.Range("E20").Value = "  +3.86%  "

$ws.Range("D21").Value = "2.096.19"
$ws.Range("E21").Value = "  -0.88%  "

$ws.Range("D22").Value = "1.002"
$ws.Range("E22").Value = "  +0.16%  "

$ws.Range("D23").Value = "5.212"
$ws.Range("E23").Value = "  -4.91%  "

$ws.Range("D24").Value = "6.127"
$ws.Range("E24").Value = "  -0.26%  "

$ws.Range("D25").Value = "166.60"
$ws.Range("E25").Value = "  +1.01%  "

$ws.Range("D26").Value = "9.151"
$ws.Range("E26").Value = "  -0.29%  "

$ws.Range("D27").Value = "18.78"
$ws.Range("E27").Value = "  -1.69%  "

$ws.Range("D28").Value = "1.923"
$ws.Range("E28").Value = "  -0.15%  "

$ws.Range("E29").Value = "  +1.50%  "

$ws.Range("D30").Value = "0.09839"
$ws.Range("E30").Value = "  +1.99%  "

$ws.Range("D31").Value = "1.464"
$ws.Range("E31").Value = "  -0.87%  "

$ws.Range("E32").Value = "  -3.05%  "

$ws.Range("D33").Value = "3.989"
$ws.Range("E33").Value = "  -2.54%  "

$ws.Range("D34").Value = "0.04672"
$ws.Range("E34").Value = "  -0.57%  "

$ws.Range("D35").Value = "1.115"
$ws.Range("E35").Value = "  -0.04%  "

$ws.Range("D36").Value = "0.6952"
$ws.Range("E36").Value = "  -1.42%  "

$ws.Range("D37").Value = "2.712"
$ws.Range("E37").Value = "  -0.58%  "

$ws.Range("D38").Value = "0.01861"
$ws.Range("E38").Value = "  +0.70%  "

$ws.Range("D39").Value = "2.601"
$ws.Range("E39").Value = "  +2.61%  "

$ws.Range("D40").Value = "6.303"
$ws.Range("E40").Value = "  +1.00%  "

$ws.Range("D41").Value = "73.11"
$ws.Range("E41").Value = "  +0.04%  "

$ws.Range("D42").Value = "1.925"
$ws.Range("E42").Value = "  -0.80%  "

$ws.Range("D44").Value = "0.8361"
$ws.Range("E44").Value = "  -0.97%  "

$ws.Range("D45").Value = "103.09"
$ws.Range("E45").Value = "  -0.71%  "

$ws.Range("D46").Value = "0.4110"
$ws.Range("E46").Value = "  -1.33%  "

$ws.Range("D47").Value = "937.77"
$ws.Range("E47").Value = "  +0.27%  "

$ws.Range("D48").Value = "9.111"
$ws.Range("E48").Value = "  -1.08%  "

$ws.Range("D49").Value = "6.954"
$ws.Range("E49").Value = "  -3.07%  "

$ws.Range("D50").Value = "33.69"
$ws.Range("E50").Value = "  -1.04%  "

$ws.Range("E51").Value = "  +0.39%  "

# Restore the original (default/general) cell formatting so no stray
# number-format styles are introduced on the updated cells.
$priceRange.ClearFormats()
